# "popravio prezentaciju i primere"
#
# Two content fixes in the "Nizovi u Javi" deck:
#
# 1. Slide 19 ("java.util.Arrays" methods overview), Rectangle 5: the
#    static method name was wrong - "copyFrom()" should read "copyOf()"
#    (the real java.util.Arrays method). We retype the "From" part of the
#    run into "Of".
#
# 2. Slide 20 (varargs slide), Rectangle 5, bullet about the vararg
#    parameter needing to be last: the sentence was split across two
#    runs ("... у списку " + "параметара"); join them back into one
#    continuous sentence/run: "... у списку параметара".

$p = $ppt.ActivePresentation

# --- Fix 1: slide 19, "copyFrom()" -> "copyOf()" -----------------------
$s19 = $p.Slides.Item(19)
$sh19 = $s19.Shapes.Item(1)
$tr19 = $sh19.TextFrame.TextRange

$full19 = $tr19.Text
$pos19 = $full19.IndexOf("copyFrom")
if ($pos19 -ge 0) {
    # "copyFrom" -> characters 1-4 are "copy", 5-8 are "From"; replace
    # just the "From" part (in place) with "Of".
    $fromRange = $tr19.Characters($pos19 + 1 + 4, 4)
    $fromRange.Text = "Of"
}

# --- Fix 2: slide 20, merge the split "...списку " / "параметара" run --
$s20 = $p.Slides.Item(20)
$sh20 = $s20.Shapes.Item(1)
$tr20 = $sh20.TextFrame.TextRange

$firstPart = "Овако дефинисан параметар функције мора да буде последњи у списку "
$secondPart = "параметара"

$full20 = $tr20.Text
$pos20 = $full20.IndexOf($firstPart)
if ($pos20 -ge 0) {
    $combined = $tr20.Characters($pos20 + 1, $firstPart.Length + $secondPart.Length)
    $combined.Text = $firstPart + $secondPart
}
